$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order of names (row 7 is a new row, status stays "Absent" for all)
$names = @("JackNickelson", "ShaquilleO'Neal", "EloneMusk", "Drake", "Rihanna", "ArianaGrande")

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = "Absent"
}
